# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (so it lands
#    before "总计"), styled like the existing "2021-Q4" sheet, and fill it
#    with the four fund rows for the new quarter.
# 2. Insert a new first data row into "总计" for "2022-Q1" (pushing the
#    existing "2021-Q4" row down), styled like the existing summary row.

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Create the "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# Headers
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Numeric-looking columns (B,D,E,F,G) must stay TEXT, like the source data
# (column B holds fund codes such as "009956" - leading zeros matter).
$wsQ1.Range("B2:B5").NumberFormat = "@"
$wsQ1.Range("D2:G5").NumberFormat = "@"

# Row 2
$wsQ1.Range("A2").Value = 0
$wsQ1.Range("B2").Value = "870009"
$wsQ1.Range("C2").Value = "广发资管平衡精选一年持有混合A"
$wsQ1.Range("D2").Value = "11.34"
$wsQ1.Range("E2").Value = "94.29"
$wsQ1.Range("F2").Value = "4.27"
$wsQ1.Range("G2").Value = "0.4842"
$wsQ1.Range("H2").Value = 8

# Row 3
$wsQ1.Range("A3").Value = 1
$wsQ1.Range("B3").Value = "872019"
$wsQ1.Range("C3").Value = "广发资管平衡精选一年持有混合C"
$wsQ1.Range("D3").Value = "1.54"
$wsQ1.Range("E3").Value = "94.29"
$wsQ1.Range("F3").Value = "4.27"
$wsQ1.Range("G3").Value = "0.0658"
$wsQ1.Range("H3").Value = 8

# Row 4
$wsQ1.Range("A4").Value = 2
$wsQ1.Range("B4").Value = "009956"
$wsQ1.Range("C4").Value = "广发恒誉混合A"
$wsQ1.Range("D4").Value = "4.94"
$wsQ1.Range("E4").Value = "21.40"
$wsQ1.Range("F4").Value = "0.93"
$wsQ1.Range("G4").Value = "0.0459"
$wsQ1.Range("H4").Value = 5

# Row 5
$wsQ1.Range("A5").Value = 3
$wsQ1.Range("B5").Value = "009957"
$wsQ1.Range("C5").Value = "广发恒誉混合C"
$wsQ1.Range("D5").Value = "0.10"
$wsQ1.Range("E5").Value = "21.40"
$wsQ1.Range("F5").Value = "0.93"
$wsQ1.Range("G5").Value = "0.0009"
$wsQ1.Range("H5").Value = 5

# Bring over header (row1) and id-column (col A) formatting from 2021-Q4.
# Done LAST (after all values are written) since pasting formats onto a
# cell and only then writing its value can clear the applied style.
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

$wsQ4.Range("A2:A5").Copy()
$wsQ1.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Update "总计": add a new row for 2022-Q1 above the 2021-Q4 row
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Move the existing "2021-Q4" row down to row 3 (it becomes the second
# entry, so its index column becomes 1), then write the new row 2.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 10
$wsTotal.Range("D3").Value = 1.54

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.6

# Re-apply the id-column style (col A) onto the (now two) data rows -
# done last, same reasoning as above.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
